{"js": "// Insert three new paragraphs right before the final (empty) trailing\n// paragraph: two plain body paragraphs, then a Heading 1 paragraph titled\n// \"The Outliner\".\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// The document ends with an empty trailing paragraph; anchor the new\n// content right before it so it lands just after the \"Write Up\" title.\nconst anchor = paragraphs.items[paragraphs.items.length - 1];\n\nconst para1 = anchor.insertParagraph(\n  \"You might think that this Outliner panel is just a simple list of what you can see in the viewport, but it is a lot more then just this. The outliner can allow something to be seen or even selected. You can organize and manage different elements in the outliner and you are even able to enter different modes from here, like selecting an item and going into the edit workspace.\",\n  Word.InsertLocation.before\n);\n\nconst para2 = para1.insertParagraph(\n  \"So, if you would like to learn a bit more about what you are actually do within this panel, then please join us for our brand-new panel entitled:\",\n  Word.InsertLocation.after\n);\n\nconst para3 = para2.insertParagraph(\"The Outliner\", Word.InsertLocation.after);\npara3.styleBuiltIn = Word.Style.heading1;\n\nawait context.sync();\n", "ps1": "# Insert three new paragraphs right before the final (empty) trailing\n# paragraph: two plain body paragraphs, then a Heading 1 paragraph titled\n# \"The Outliner\".\n$d = $word.ActiveDocument\n\n$lastPara = $d.Paragraphs.Last\n\n# Create three blank paragraphs just ahead of the trailing empty paragraph\n# (so the new content lands right after the \"Write Up\" title).\n$lastPara.Range.InsertParagraphBefore()\n$lastPara.Range.InsertParagraphBefore()\n$lastPara.Range.InsertParagraphBefore()\n\n$para1 = $d.Paragraphs.Item(2)\n$para1.Range.Text = \"You might think that this Outliner panel is just a simple list of what you can see in the viewport, but it is a lot more then just this. The outliner can allow something to be seen or even selected. You can organize and manage different elements in the outliner and you are even able to enter different modes from here, like selecting an item and going into the edit workspace.\"\n\n$para2 = $d.Paragraphs.Item(3)\n$para2.Range.Text = \"So, if you would like to learn a bit more about what you are actually do within this panel, then please join us for our brand-new panel entitled:\"\n\n$para3 = $d.Paragraphs.Item(4)\n$para3.Range.Text = \"The Outliner\"\n$para3.Style = \"Heading 1\"\n"}
